$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# --- Content edit -------------------------------------------------------
# The "CasesTab" query in B2 dropped its trailing `Cohort` column
# (coalesce(co.cohort_description, '') AS `Cohort`) from the RETURN list.
$queryText = $ws.Range("B2").Text
$cohortFragment = ",`n        coalesce(co.cohort_description, '') AS ``Cohort``"
if ($queryText.Contains($cohortFragment)) {
    $queryText = $queryText.Replace($cohortFragment, "")
} else {
    # Fall back to a tolerant regex trim in case of whitespace/newline drift.
    $queryText = [System.Text.RegularExpressions.Regex]::Replace(
        $queryText,
        ",\s*coalesce\(co\.cohort_description,\s*''\)\s*AS\s*``Cohort``\s*$",
        "")
}
$ws.Range("B2").Value = $queryText

# --- Row heights ----------------------------------------------------------
# Rows with the wrapped query text were reflowed to slightly different
# heights in the resave.
$ws.Rows.Item(2).RowHeight = 300
$ws.Rows.Item(3).RowHeight = 285
$ws.Rows.Item(4).RowHeight = 270
